$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-02-25 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-26 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("95×98=", $true, $false, $false, $false, $false, $true, 1, $false, "28×75=", 2) | Out-Null
$d.Content.Find.Execute("89×83=", $true, $false, $false, $false, $false, $true, 1, $false, "80×39=", 2) | Out-Null
$d.Content.Find.Execute("50×20=", $true, $false, $false, $false, $false, $true, 1, $false, "84×70=", 2) | Out-Null
$d.Content.Find.Execute("97×93=", $true, $false, $false, $false, $false, $true, 1, $false, "50×78=", 2) | Out-Null
$d.Content.Find.Execute("55×60=", $true, $false, $false, $false, $false, $true, 1, $false, "98×32=", 2) | Out-Null
$d.Content.Find.Execute("96×99=", $true, $false, $false, $false, $false, $true, 1, $false, "39×91=", 2) | Out-Null
$d.Content.Find.Execute("36×97=", $true, $false, $false, $false, $false, $true, 1, $false, "27×32=", 2) | Out-Null
$d.Content.Find.Execute("68×68=", $true, $false, $false, $false, $false, $true, 1, $false, "74×89=", 2) | Out-Null
$d.Content.Find.Execute("72×43=", $true, $false, $false, $false, $false, $true, 1, $false, "63×14=", 2) | Out-Null
$d.Content.Find.Execute("65×70=", $true, $false, $false, $false, $false, $true, 1, $false, "20×64=", 2) | Out-Null
$d.Content.Find.Execute("96×43=", $true, $false, $false, $false, $false, $true, 1, $false, "69×54=", 2) | Out-Null
$d.Content.Find.Execute("45×64=", $true, $false, $false, $false, $false, $true, 1, $false, "82×20=", 2) | Out-Null
$d.Content.Find.Execute("50×37=", $true, $false, $false, $false, $false, $true, 1, $false, "98×94=", 2) | Out-Null
$d.Content.Find.Execute("95×33=", $true, $false, $false, $false, $false, $true, 1, $false, "57×82=", 2) | Out-Null
$d.Content.Find.Execute("39×71=", $true, $false, $false, $false, $false, $true, 1, $false, "36×46=", 2) | Out-Null
$d.Content.Find.Execute("49×95=", $true, $false, $false, $false, $false, $true, 1, $false, "35×64=", 2) | Out-Null
$d.Content.Find.Execute("86×40=", $true, $false, $false, $false, $false, $true, 1, $false, "78×23=", 2) | Out-Null
$d.Content.Find.Execute("38×26=", $true, $false, $false, $false, $false, $true, 1, $false, "49×92=", 2) | Out-Null
$d.Content.Find.Execute("41×78=", $true, $false, $false, $false, $false, $true, 1, $false, "64×26=", 2) | Out-Null
$d.Content.Find.Execute("99×57=", $true, $false, $false, $false, $false, $true, 1, $false, "29×99=", 2) | Out-Null
$d.Content.Find.Execute("45×91=", $true, $false, $false, $false, $false, $true, 1, $false, "98×36=", 2) | Out-Null
$d.Content.Find.Execute("18×31=", $true, $false, $false, $false, $false, $true, 1, $false, "16×83=", 2) | Out-Null
$d.Content.Find.Execute("14×61=", $true, $false, $false, $false, $false, $true, 1, $false, "11×53=", 2) | Out-Null
$d.Content.Find.Execute("93×69=", $true, $false, $false, $false, $false, $true, 1, $false, "47×92=", 2) | Out-Null
$d.Content.Find.Execute("54×89=", $true, $false, $false, $false, $false, $true, 1, $false, "95×55=", 2) | Out-Null
